$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by duplicating "2022-Q3" (so it
#    inherits the same column layout/styles), placing it right before
#    the existing "2022-Q3" tab, then overwrite its data with the
#    Q4 figures.
# ------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)

$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# The source sheet has 3 data rows (rows 2-4); Q4 only needs 1 (row 2).
$q4.Rows.Item(4).Delete()
$q4.Rows.Item(3).Delete()

# Write the Q4 fund row as text (matches the inline-string typing used
# for the analogous columns on every other quarter sheet) and restore
# the default (unstyled) look of the template's data-row cells.
$q4.Range("B2:G2").NumberFormat = "@"
$q4.Range("B2").Value = "003397"
$q4.Range("C2").Value = "银华体育文化灵活配置混合"
$q4.Range("D2").Value = "0.53"
$q4.Range("E2").Value = "87.74"
$q4.Range("F2").Value = "3.88"
$q4.Range("G2").Value = "0.0206"
$q4.Range("H2").Value = 5
$q4.Range("B2:G2").Style = "Normal"

# ------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q4
#    above the existing data and renumber the index column.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()

# Re-apply the data-row formatting (copied from row 3, which still has
# it) to the freshly inserted row 2 so it matches the rest of the
# table instead of inheriting the header row's style.
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.02

# Renumber the remaining rows' index column (0-based, sequential).
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5

# Restore "总计" as the active tab (adding/copying sheets along the way
# shifts Excel's notion of the active sheet to whichever was touched
# last) so the book view matches the original, unmodified state.
$total.Activate()

Write-Output "2022-Q4 sheet added"
